$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$old = "What are the outcomes of these incidents? Which states with the highest number of killed and injured due to gun related violence?"
$placeholder = "XPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERXPLACEHOLDERX"
$new = "What are the outcomes of these incidents? Which states have the highest number of killed and injured due to gun related violence?"

$full = $tr.Text
$full1 = $full.Replace($old, $placeholder)
$tr.Text = $full1
Write-Output "step1=[$($shp.TextFrame.TextRange.Text)]"

$full2 = $shp.TextFrame.TextRange.Text
$full3 = $full2.Replace($placeholder, $new)
$shp.TextFrame.TextRange.Text = $full3
Write-Output "step2=[$($shp.TextFrame.TextRange.Text)]"
